$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40/41: VeChain and dogwifhat swap positions, with updated price/volume.
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0395"
$ws.Range("E40").Value = "  -7.69%  "

$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.87"
$ws.Range("E41").Value = "  -14.06%  "

# Price / volume refresh for every other row.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.366.09"
$ws.Range("E2").Value = "  -3.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.153.86"
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.91"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  -6.69%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.147.58"
$ws.Range("E8").Value = "  -3.00%  "
$ws.Range("E9").Value = "  -4.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("E10").Value = "  -6.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.51"
$ws.Range("E11").Value = "  -5.53%  "
$ws.Range("E12").Value = "  -5.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  -8.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.08"
$ws.Range("E14").Value = "  -8.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.669.95"
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.370.93"
$ws.Range("E16").Value = "  -3.51%  "
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.149.24"
$ws.Range("E18").Value = "  -2.90%  "
$ws.Range("E19").Value = "  -5.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.65"
$ws.Range("E20").Value = "  -5.89%  "
$ws.Range("E21").Value = "  -4.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.709"
$ws.Range("E22").Value = "  -5.43%  "
$ws.Range("E23").Value = "  -4.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.75"
$ws.Range("E24").Value = "  -6.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.50"
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.87"
$ws.Range("E27").Value = "  -4.81%  "
$ws.Range("E28").Value = "  -7.68%  "
$ws.Range("E29").Value = "  -7.10%  "
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.114"
$ws.Range("E31").Value = "  -33.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.76"
$ws.Range("E32").Value = "  -3.41%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.21"
$ws.Range("E35").Value = "  -4.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.36"
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("E37").Value = "  -6.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0711"
$ws.Range("E38").Value = "  -11.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "446.32"
$ws.Range("E39").Value = "  -9.40%  "
$ws.Range("E42").Value = "  -7.76%  "
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.834.08"
$ws.Range("E44").Value = "  -3.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.267"
$ws.Range("E45").Value = "  -9.36%  "
$ws.Range("E46").Value = "  -8.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.41"
$ws.Range("E47").Value = "  -7.01%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("E49").Value = "  -4.24%  "
$ws.Range("E50").Value = "  -4.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "117.99"
$ws.Range("E51").Value = "  -2.97%  "
